# Pandas cheat sheet: fix "Pd.notnull" -> "pd.notnull" in the
# "Logic in Python (and pandas)" table on slide 1 (row 5, column 3).
#
# The cell's paragraph holds four separate runs ("Pd.notnull", "(", "obj",
# ")"); only the text of the first run changes. Plain TextRange.Text / Find /
# Characters assignment on a multi-run table-cell paragraph rewrites the
# first touched run with the *whole* new string and then leaves the other
# original runs behind (duplicating their text), so instead we address the
# individual run directly via TextFrame2's run collection and only touch
# that run's Text - this leaves the sibling runs (and their own formatting)
# completely untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Table 13") {
        $tbl = $sh.Table
        $cellShape = $tbl.Cell(5, 3).Shape
        $run = $cellShape.TextFrame2.TextRange.Runs(1, 1)
        if ($run.Text.StartsWith("Pd.notnull")) {
            $run.Text = "pd.notnull"
        }
    }
}
